# Apply the "New file structure, setup.py, finished new weighted histogram
# matching method" update to the Experiments sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# --- Update the note text on row 44 (L44) ---------------------------------
$ws.Range("L44").Value2 = "Does not match what is published: UPDATE: Because github model is close but not quite as good as what went in the paper."

# Row 44 grew taller to fit the longer note text.
$ws.Rows.Item(44).RowHeight = 31

# --- New row 46: DenseDepth (Wasserstein histogram matching) label --------
# (Entered first so the shared-string table order matches the source file.)
$ws.Range("A46").Value2 = "DenseDepth (Wasserstein histogram matching)"

# --- New row 41: DORN (Wasserstein histogram matching) results ------------
$ws.Range("A41").Value2 = "DORN (Wasserstein histogram matching)"
$ws.Range("B41").Value2 = 0.847427449419342
$ws.Range("C41").Value2 = 0.95332383895321304
$ws.Range("D41").Value2 = 0.982672920285379
$ws.Range("F41").Value2 = 0.49932645306856899
$ws.Range("G41").Value2 = 0.117189220622728
$ws.Range("I41").Value2 = 0.053489108434636203

# Match number formatting of the neighboring numeric cells (style index 2 / "0.000").
$ws.Range("B41:D41").NumberFormat = "0.000"
$ws.Range("F41:G41").NumberFormat = "0.000"
$ws.Range("I41").NumberFormat = "0.000"

# --- Update the sheet view's scroll position / active cell selection ------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 30
$ws.Range("J45").Select()
